$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Range("H137").Value = 2279.44
$ws.Range("I137").Value = 1989.05
$ws.Range("J137").Value = 3441
$ws.Range("K137").Value = 5967.15
$ws.Range("L137").Value = 10323
$ws.Range("M137").Value = -3417.15
$ws.Range("N137").Value = -15423

# Row 138
$ws.Range("H138").Value = 4190.47
$ws.Range("I138").Value = 3337.9333
$ws.Range("J138").Value = 4888
$ws.Range("K138").Value = 10013.7999
$ws.Range("L138").Value = 14664
$ws.Range("M138").Value = -4873.7999
$ws.Range("N138").Value = -24944

$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 32
$ws.Range("H32").Value = 26821.316
$ws.Range("I32").Value = 21272.59
$ws.Range("J32").Value = 104503.5
$ws.Range("K32").Value = 21272.59
$ws.Range("L32").Value = 104503.5
$ws.Range("M32").Value = -20985.59
$ws.Range("N32").Value = -105077.5

# Row 61
$ws.Range("H61").Value = 2452.6572
$ws.Range("I61").Value = 2078.3333
$ws.Range("J61").Value = 3716
$ws.Range("K61").Value = 2078.3333
$ws.Range("L61").Value = 3716
$ws.Range("M61").Value = -1866.3333
$ws.Range("N61").Value = -4140

# Row 132
$ws.Range("H132").Value = 2115.1396
$ws.Range("I132").Value = 1733.2778
$ws.Range("J132").Value = 4079
$ws.Range("K132").Value = 5199.8334
$ws.Range("L132").Value = 12237
$ws.Range("M132").Value = -2669.8334
$ws.Range("N132").Value = -17297

# Row 136
$ws.Range("H136").Value = 2452.6572
$ws.Range("I136").Value = 2078.3333
$ws.Range("J136").Value = 3716
$ws.Range("K136").Value = 6234.999899999999
$ws.Range("L136").Value = 11148
$ws.Range("M136").Value = -3684.999899999999
$ws.Range("N136").Value = -16248

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 2989.3157
$ws.Range("I22").Value = 3915.2144
$ws.Range("J22").Value = 396.8
$ws.Range("K22").Value = 3915.2144
$ws.Range("L22").Value = 396.8
$ws.Range("M22").Value = -3742.2144
$ws.Range("N22").Value = -742.8

# Row 107
$ws.Range("H107").Value = 21744.346
$ws.Range("I107").Value = 36157.4
$ws.Range("J107").Value = 2090.182
$ws.Range("K107").Value = 36157.4
$ws.Range("L107").Value = 2090.182
$ws.Range("M107").Value = -34237.4
$ws.Range("N107").Value = -5930.182

# Row 132
$ws.Range("H132").Value = 75094.28999999999
$ws.Range("J132").Value = 75094.28999999999
$ws.Range("L132").Value = 75094.28999999999
$ws.Range("N132").Value = -85214.28999999999

# Row 140
$ws.Range("H140").Value = 63187.145
$ws.Range("J140").Value = 63187.145
$ws.Range("L140").Value = 63187.145
$ws.Range("N140").Value = -73547.14499999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5309.2
$ws.Range("I31").Value = 4001.75
$ws.Range("J31").Value = 6180.8335
$ws.Range("K31").Value = 4001.75
$ws.Range("L31").Value = 6180.8335
$ws.Range("M31").Value = -3706.75
$ws.Range("N31").Value = -6770.8335

# Row 34
$ws.Range("H34").Value = 5309.2
$ws.Range("I34").Value = 4001.75
$ws.Range("J34").Value = 6180.8335
$ws.Range("K34").Value = 4001.75
$ws.Range("L34").Value = 6180.8335
$ws.Range("M34").Value = -3799.75
$ws.Range("N34").Value = -6584.8335

# Row 58
$ws.Range("H58").Value = 1290.05
$ws.Range("I58").Value = 1396.0968
$ws.Range("J58").Value = 924.7778
$ws.Range("K58").Value = 1396.0968
$ws.Range("L58").Value = 924.7778
$ws.Range("M58").Value = -1193.0968
$ws.Range("N58").Value = -1330.7778

# Row 107
$ws.Range("H107").Value = 885.75
$ws.Range("I107").Value = 937
$ws.Range("J107").Value = 868.6667
$ws.Range("K107").Value = 937
$ws.Range("L107").Value = 868.6667
$ws.Range("M107").Value = 983
$ws.Range("N107").Value = -4708.6667

# Row 136
$ws.Range("H136").Value = 1290.05
$ws.Range("I136").Value = 1396.0968
$ws.Range("J136").Value = 924.7778
$ws.Range("K136").Value = 4188.2904
$ws.Range("L136").Value = 2774.3334
$ws.Range("M136").Value = -1638.2904
$ws.Range("N136").Value = -7874.3334

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# Row 141
$ws.Range("H141").Value = 40447.39
$ws.Range("J141").Value = 40473.65
$ws.Range("L141").Value = 40473.65
$ws.Range("N141").Value = -50833.65

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1730.6857
$ws.Range("I5").Value = 1384
$ws.Range("J5").Value = 2142.375
$ws.Range("K5").Value = 4152
$ws.Range("L5").Value = 6427.125
$ws.Range("M5").Value = -4040
$ws.Range("N5").Value = -6651.125

# Row 23
$ws.Range("H23").Value = 399.8889
$ws.Range("I23").Value = 300
$ws.Range("J23").Value = 428.42856
$ws.Range("K23").Value = 900
$ws.Range("L23").Value = 1285.28568
$ws.Range("M23").Value = -665
$ws.Range("N23").Value = -1755.28568

# Row 132
$ws.Range("H132").Value = 1837.5
$ws.Range("I132").Value = 1428.909
$ws.Range("K132").Value = 12860.181
$ws.Range("M132").Value = -10330.181

# Row 135
$ws.Range("H135").Value = 1730.6857
$ws.Range("I135").Value = 1384
$ws.Range("J135").Value = 2142.375
$ws.Range("K135").Value = 12456
$ws.Range("L135").Value = 19281.375
$ws.Range("M135").Value = -9921
$ws.Range("N135").Value = -24351.375

$ws = $wb.Worksheets.Item("GSM")
# Row 36
$ws.Range("H36").Value = 6667.067
$ws.Range("I36").Value = 2997
$ws.Range("J36").Value = 8001.636
$ws.Range("K36").Value = 2997
$ws.Range("L36").Value = 8001.636
$ws.Range("M36").Value = -2512
$ws.Range("N36").Value = -8971.636

# Row 70
$ws.Range("H70").Value = 6490.85
$ws.Range("I70").Value = 4961.6
$ws.Range("J70").Value = 7000.6
$ws.Range("K70").Value = 4961.6
$ws.Range("L70").Value = 7000.6
$ws.Range("M70").Value = -4691.6
$ws.Range("N70").Value = -7540.6

# Row 73
$ws.Range("H73").Value = 6490.85
$ws.Range("I73").Value = 4961.6
$ws.Range("J73").Value = 7000.6
$ws.Range("K73").Value = 4961.6
$ws.Range("L73").Value = 7000.6
$ws.Range("M73").Value = -4025.6
$ws.Range("N73").Value = -8872.6

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 17612.625
$ws.Range("I136").Value = 24379.2
$ws.Range("K136").Value = 73137.60000000001
$ws.Range("M136").Value = -70587.60000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 7759
$ws.Range("I136").Value = 9198.75
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 27596.25
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -4549.250100000001
$ws.Range("N136").Value = -11100
